# Generate Report for Handback
#
# For both locale sheets ("zh-cn" and "de-de") this:
#   - updates the Status column (B) from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two real rows,
#   - fills in the "Latest Target File" (E) and "Latest Handback File" (F)
#     columns (previously empty) with hyperlinked filenames mirroring the
#     Source File Name (A) / Latest Handoff File (C) columns,
#   - stamps the "Latest Handback DateTime" (G) column with the actual
#     handback timestamp (previously the zero-date sentinel).

$wb = $excel.ActiveWorkbook

# The "Overview" summary sheet mirrors the same Status text for each
# locale/file pair (it shares the same underlying string), so it also
# flips from "Ready for handoff" to "Handed back: in sync with en-US".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$localeSheets = @(
    @{ SheetName = "zh-cn"; HandbackTime2 = "2016-03-07 02:19:12"; HandbackTime3 = "2016-03-07 02:19:12" },
    @{ SheetName = "de-de"; HandbackTime2 = "2016-03-07 02:19:30"; HandbackTime3 = "2016-03-07 02:19:30" }
)

foreach ($cfg in $localeSheets) {
    $ws = $wb.Worksheets.Item($cfg.SheetName)

    # Collect the existing hyperlink addresses (in sheet order: A2, C2, A3, C3, A4)
    # so the new Target File / Handback File links can reuse them.
    $addrs = @()
    foreach ($h in $ws.Hyperlinks) {
        $addrs += $h.Address
    }
    $addrA2 = $addrs[0]
    $addrC2 = $addrs[1]
    $addrA3 = $addrs[2]
    $addrC3 = $addrs[3]

    $fileA2 = $ws.Range("A2").Text
    $fileC2 = $ws.Range("C2").Text
    $fileA3 = $ws.Range("A3").Text
    $fileC3 = $ws.Range("C3").Text

    # Status -> handed back, now in sync with en-US
    $ws.Range("B2").Value = "Handed back: in sync with en-US"
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    # Latest Target File / Latest Handback File for row 2, mirroring row 2's
    # source file / handoff file (same underlying document, now targeted +
    # handed back).
    $ws.Hyperlinks.Add($ws.Range("E2"), $addrA2, "", "", $fileA2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrC2, "", "", $fileC2)

    # Latest Target File / Latest Handback File for row 3.
    $ws.Hyperlinks.Add($ws.Range("E3"), $addrA3, "", "", $fileA3)
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrC3, "", "", $fileC3)

    # Latest Handback DateTime now reflects the real handback timestamps.
    $ws.Range("G2").Value = $cfg.HandbackTime2
    $ws.Range("G3").Value = $cfg.HandbackTime3
}
